# Remove the "operator_key" entry from the fact_electricity column (A5) of
# the data dictionary sheet. Only column A is affected: the values in
# A6:A28 shift up one row to A5:A27, leaving A28 empty. Columns C, E and G
# (and their row positions) are left completely untouched, so a plain
# "delete cell, shift whole row up" cannot be used here - each column's
# values are rewritten individually instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 28

# Read the existing column-A values (row 6 downward) before overwriting
# anything, then write them back starting one row higher - this is the
# "delete A5, shift A6:A28 up" operation restricted to column A only.
# (Value2 is used for the read because it reliably returns the scalar
# contents; Value is used for the write.)
for ($r = 6; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r - 1, 1).Value = $ws.Cells.Item($r, 1).Value2
}

# The bottom row no longer has a column-A entry once everything shifted up.
$ws.Cells.Item($lastRow, 1).ClearContents()

# Reset the saved selection to the default top-left cell.
$ws.Range("A1").Select() | Out-Null
